$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8-11: Team member names (replacing "Member 1".."Member 4"); Row 12 cleared (was "Member 5")
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("A12").Value = ""

$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100

# Row 3: Date value
$ws.Range("B3").Value = [DateTime]"2020-10-01"

# Row 4: Team Name
$ws.Range("B4").Value = "Limette"

# Row 5: Total Number of Team Members
$ws.Range("B5").Value = 4

# Row 19-21: Tasks
$ws.Range("A19").Value = "Interviews"
$ws.Range("A20").Value = "Affinity Diagram"
$ws.Range("B19").Value = "Brainstorming"
$ws.Range("B20").Value = "Personas"
$ws.Range("B21").Value = "Presentation"

# Row 18 re-wraps at a shorter height under the new Excel build's text metrics
$ws.Rows.Item(18).RowHeight = 37.15

# View adjustments
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C20").Select()
